$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H) — copy the format from the existing "sum" header (G1)
# so it picks up the same bold/centered/bordered style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column H2:H6 — all zeros (no special style, like the other data cells)
$ws.Range("H2:H6").Value = 0
